$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Luis Escudero's name to a linked markdown reference
$ws.Range("A31").Value = "[Luis Escudero](https://www.gob.pe/institucion/imarpe/funcionarios/30904-luis-orlando-escudero-herrera)"

# Correct the malformed CEOAS - OSU markdown link
$ws.Range("C34").Value = "[CEOAS - OSU](https://ceoas.oregonstate.edu/)"

# Replace the plain affiliation text with the UNICAMP markdown link
$ws.Range("C43").Value = "[UNICAMP](https://www.ib.unicamp.br/)"

# Update the view's scroll position / selection to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("A35").Select()
